$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    2.356563929135413,
    2.239645347542925,
    2.309403061943285,
    2.324704126691745,
    2.393542322469567,
    2.305544188113036,
    2.144322488314476,
    2.164041425356085,
    1.893541085346338,
    1.898704885206921,
    1.88122482540276,
    1.861186286886721,
    2.060335824929814,
    1.993352207562557,
    2.188624724064947,
    1.879312442730859,
    1.833104773380641,
    1.831963517995406,
    1.756092968114026,
    1.485188413353329
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Remove rows 22-32 (previously held Path values 20..34)
$ws.Range("A22:B32").EntireRow.Delete() | Out-Null
